# Add 2022-Q3 data: insert a new sheet with the latest quarter's fund
# holdings and refresh the "总计" (summary) sheet with the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: shift every existing quarter down
#    one slot and put the brand-new 2022-Q3 figures in row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 7.52

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 5
$summary.Range("D3").Value = 7.29

$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 3
$summary.Range("D4").Value = 6.44

$summary.Range("B5").Value = "2021-Q3"
$summary.Range("C5").Value = 3
$summary.Range("D5").Value = 3.6

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q2"
$summary.Range("C6").Value = 3
$summary.Range("D6").Value = 2.79

# Copy row-2's formatting down onto the freshly-used row 6 (same look as
# the other data rows: centred bold "A" cell style).
$summary.Range("A2").Copy()
$summary.Range("A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Insert a brand new "2022-Q3" sheet right after "总计" — every other
#    quarter tab simply shifts one position to the right automatically.
# ---------------------------------------------------------------------
$styleSource = $wb.Worksheets.Item("2022-Q2")

$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Data rows — fund code / name / scale / position / ratio / market value
# are text in the source workbook (keeps trailing zeros & leading
# zeros), only the rank column ("H") is numeric.
$rows = @(
    @(0, "012349", "天弘恒生科技指数（QDII）C", "33.57", "92.84", "8.32", "2.7930", 3),
    @(1, "012348", "天弘恒生科技指数（QDII）A", "30.64", "92.84", "8.32", "2.5492", 3),
    @(2, "968029", "恒生指数基金M类人民币（对冲）份额", "27.03", "99.07", "6.93", "1.8732", 5),
    @(3, "009225", "天弘中证中美互联网指数（QDII）A", "1.20", "94.98", "9.18", "0.1102", 3),
    @(4, "002379", "工银瑞信香港中小盘股票（QDII）人民币", "1.58", "78.58", "4.38", "0.0692", 4),
    @(5, "002380", "工银瑞信香港中小盘股票（QDII）美元", "1.58", "78.58", "4.38", "0.0692", 4),
    @(6, "009226", "天弘中证中美互联网指数（QDII）C", "0.60", "94.98", "9.18", "0.0551", 3)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
}

# Match the look & feel of the other quarter sheets (bold centred header
# row + bold centred "A" column).
$styleSource.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$styleSource.Range("A2").Copy()
$q3.Range("A2:A8").PasteSpecial(-4122)

# Keep the original active/selected tab ("2021-Q2", now the last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
$lastSheet.Range("A1").Select()

Write-Host "2022-Q3 sheet added; summary refreshed."
